# Edit slide 1 / shape id 2050 ("Rectangle 2", the ctrTitle placeholder):
#   - resize/reposition the title box
#   - shrink the title font from 54pt to 44pt
#   - replace "Baze de date" with the longer, multi-run title
#     "Introducere în analiza datelor de mari dimensiuni"
#     (kept as separate runs per word/space, matching the original
#     per-word run layout used by the author's spell-checked edit)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

# --- reposition / resize the placeholder (EMU 762000/2057399/8458200/2514601) ---
$shp.Left = 60
$shp.Top = 161.99992
$shp.Width = 666
$shp.Height = 198.0001

# --- replace the text, then size the whole run, then re-touch each
#     word/space sub-range so the engine keeps them as distinct runs
#     (mirrors the separate <a:r> runs seen in the final file) ---
$tr = $shp.TextFrame.TextRange
$tr.Text = "Introducere în analiza datelor de mari dimensiuni"
$tr.Font.Size = 44

$tr.Characters(1, 11).Font.Size = 44    # Introducere
$tr.Characters(12, 1).Font.Size = 44    # " "
$tr.Characters(13, 2).Font.Size = 44    # în
$tr.Characters(15, 1).Font.Size = 44    # " "
$tr.Characters(16, 7).Font.Size = 44    # analiza
$tr.Characters(23, 1).Font.Size = 44    # " "
$tr.Characters(24, 7).Font.Size = 44    # datelor
$tr.Characters(31, 4).Font.Size = 44    # " de "
$tr.Characters(35, 4).Font.Size = 44    # mari
$tr.Characters(39, 1).Font.Size = 44    # " "
$tr.Characters(40, 10).Font.Size = 44   # dimensiuni
